$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.08528554791034348
$ws.Range("C2").Value = 0.4252855479103435
$ws.Range("D2").Value = 0.7552855479103435
$ws.Range("E2").Value = 0.6731635053256353
$ws.Range("F2").Value = 0.3446377147317632
$ws.Range("G2").Value = 0.6452855479103434
$ws.Range("H2").Value = 0.4212855479103435
$ws.Range("B3").Value = 0.34
$ws.Range("C3").Value = 0.6699999999999999
$ws.Range("D3").Value = 0.5878779574152918
$ws.Range("E3").Value = 0.2593521668214197
$ws.Range("F3").Value = 0.5600000000000001
$ws.Range("G3").Value = 0.336
$ws.Range("B4").Value = 0.33
$ws.Range("C4").Value = 0.2478779574152918
$ws.Range("D4").Value = -0.08064783317858029
$ws.Range("E4").Value = 0.22
$ws.Range("F4").Value = -0.003999999999999997
$ws.Range("G4").Value = 0.1225100404635037
$ws.Range("H4").Value = -0.2182525219575302
$ws.Range("I4").Value = 0.08717552522494373
$ws.Range("J4").Value = -0.1742297805489477
$ws.Range("B5").Value = -0.08212204258470818
$ws.Range("C5").Value = -0.4106478331785803
$ws.Range("D5").Value = -0.11
$ws.Range("E5").Value = -0.334
$ws.Range("F5").Value = -0.2074899595364962
$ws.Range("G5").Value = -0.5482525219575302
$ws.Range("H5").Value = -0.2428244747750563
$ws.Range("I5").Value = -0.5042297805489477
$ws.Range("B6").Value = -0.3285257905938721
$ws.Range("C6").Value = -0.02787795741529181
$ws.Range("D6").Value = -0.2518779574152918
$ws.Range("E6").Value = -0.1253679169517881
$ws.Range("F6").Value = -0.466130479372822
$ws.Range("G6").Value = -0.1607024321903481
$ws.Range("H6").Value = -0.4221077379642395
$ws.Range("B7").Value = 0.3006478331785803
$ws.Range("C7").Value = 0.0766478331785803
$ws.Range("D7").Value = 0.203157873642084
$ws.Range("E7").Value = -0.1376046887789499
$ws.Range("F7").Value = 0.167823358403524
$ws.Range("G7").Value = -0.0935819473703674
$ws.Range("B8").Value = -0.224
$ws.Range("C8").Value = -0.09748995953649625
$ws.Range("D8").Value = -0.4382525219575302
$ws.Range("E8").Value = -0.1328244747750563
$ws.Range("F8").Value = -0.3942297805489477
$ws.Range("G8").Value = -0.4194371574146135
$ws.Range("H8").Value = -0.2271788341830432
$ws.Range("I8").Value = -0.2759495356205764
$ws.Range("B9").Value = 0.1265100404635037
$ws.Range("C9").Value = -0.2142525219575302
$ws.Range("D9").Value = 0.09117552522494374
$ws.Range("E9").Value = -0.1702297805489477
$ws.Range("F9").Value = -0.1954371574146135
$ws.Range("G9").Value = -0.003178834183043253
$ws.Range("H9").Value = -0.0519495356205764
$ws.Range("B10").Value = -0.3407625624210339
$ws.Range("C10").Value = -0.03533451523856001
$ws.Range("D10").Value = -0.2967398210124514
$ws.Range("E10").Value = -0.3219471978781172
$ws.Range("F10").Value = -0.129688874646547
$ws.Range("G10").Value = -0.1784595760840801
$ws.Range("B11").Value = 0.3054280471824739
$ws.Range("C11").Value = 0.04402274140858248
$ws.Range("D11").Value = 0.01881536454291668
$ws.Range("E11").Value = 0.2110736877744869
$ws.Range("F11").Value = 0.1623029863369538
$ws.Range("B12").Value = -0.2614053057738914
$ws.Range("C12").Value = -0.2866126826395572
$ws.Range("D12").Value = -0.09435435940798698
$ws.Range("E12").Value = -0.1431250608455201
$ws.Range("B13").Value = -0.02520737686566579
$ws.Range("C13").Value = 0.1670509463659045
$ws.Range("D13").Value = 0.1182802449283713
$ws.Range("B14").Value = 0.1922583232315702
$ws.Range("C14").Value = 0.1434876217940371
$ws.Range("B15").Value = -0.04877070143753315
